# Update the MPA Adjustment "Data" sheet with the new asset numbers / subnumbers
# produced by the test automation upload run.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Main Asset Number (column D) corrections: 60000372 -> 60000390 (rows 6-10, 16-17)
$ws.Range("D6").Value = 60000390
$ws.Range("D7").Value = 60000390
$ws.Range("D8").Value = 60000390
$ws.Range("D9").Value = 60000390
$ws.Range("D10").Value = 60000390
$ws.Range("D16").Value = 60000390
$ws.Range("D17").Value = 60000390

# Asset Subnumber (column E) corrections: 270 -> 277 (rows 11-15, 18-19)
$ws.Range("E11").Value = 277
$ws.Range("E12").Value = 277
$ws.Range("E13").Value = 277
$ws.Range("E14").Value = 277
$ws.Range("E15").Value = 277
$ws.Range("E18").Value = 277
$ws.Range("E19").Value = 277

# Remaining rows each shift Main Asset Number / Asset Subnumber by +18
$ws.Range("D20").Value = 60000391
$ws.Range("E21").Value = 278
$ws.Range("D22").Value = 60000392
$ws.Range("E23").Value = 279
$ws.Range("D24").Value = 60000393
$ws.Range("E25").Value = 280
$ws.Range("D26").Value = 60000394
$ws.Range("E27").Value = 281
